$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.128.46'
$ws.Range('E2').Value = '  -0.32%  '

# Row 3
$ws.Range('D3').Value = '2.074.67'
$ws.Range('E3').Value = '  -1.12%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.06'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.675'
$ws.Range('E6').Value = '  +1.74%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.27'
$ws.Range('E7').Value = '  +9.51%  '

# Row 8
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.390'
$ws.Range('E9').Value = '  +4.00%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '61.37'
$ws.Range('E10').Value = '  -0.44%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0799'
$ws.Range('E11').Value = '  +7.45%  '

# Row 12
$ws.Range('E12').Value = '  +2.48%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '16.26'
$ws.Range('E13').Value = '  +6.08%  '

# Row 14
$ws.Range('D14').Value = '2.376.93'
$ws.Range('E14').Value = '  -1.10%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.822'
$ws.Range('E15').Value = '  -2.08%  '

# Row 16
$ws.Range('E16').Value = '  +6.69%  '

# Row 17
$ws.Range('D17').Value = '2.069.63'
$ws.Range('E17').Value = '  -1.43%  '

# Row 18
$ws.Range('D18').Value = '37.100.70'
$ws.Range('E18').Value = '  -0.38%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.70'
$ws.Range('E19').Value = '  +7.56%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '74.70'
$ws.Range('E20').Value = '  +2.93%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0930'
$ws.Range('E21').Value = '  +10.99%  '

# Row 22
$ws.Range('E22').Value = '  +4.48%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.13'
$ws.Range('E23').Value = '  -0.74%  '

# Row 24
$ws.Range('E24').Value = '  +0.02%  '

# Row 25
$ws.Range('E25').Value = '  -3.10%  '

# Row 26
$ws.Range('E26').Value = '  +12.73%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.73'
$ws.Range('E27').Value = '  -0.74%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.34'
$ws.Range('E28').Value = '  +1.13%  '

# Row 29
$ws.Range('E29').Value = '  -1.76%  '

# Row 30
$ws.Range('E30').Value = '  +2.59%  '

# Row 31
$ws.Range('E31').Value = '  +6.02%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.81'
$ws.Range('E32').Value = '  +6.31%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0634'
$ws.Range('E33').Value = '  +3.31%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.48'
$ws.Range('E34').Value = '  +8.65%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0905'
$ws.Range('E35').Value = '  -0.28%  '

# Row 36
$ws.Range('E36').Value = '  -0.07%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('E37').Value = '  -0.78%  '

# Row 38
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.78'
$ws.Range('E38').Value = '  -3.69%  '

# Row 39
$ws.Range('B39').Value = 'Cronos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.116'
$ws.Range('E39').Value = '  +26.93%  '

# Row 40
$ws.Range('E40').Value = '  +0.71%  '

# Row 41
$ws.Range('E41').Value = '  +0.48%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.76'
$ws.Range('E42').Value = '  -3.09%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.92'

# Row 45
$ws.Range('E45').Value = '  +12.22%  '

# Row 46
$ws.Range('E46').Value = '  +1.34%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.51'
$ws.Range('E47').Value = '  +12.21%  '

# Row 48
$ws.Range('E48').Value = '  +8.27%  '

# Row 49
$ws.Range('D49').Value = '1.302.56'
$ws.Range('E49').Value = '  -1.34%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('E50').Value = '  -0.99%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.94'
$ws.Range('E51').Value = '  -1.58%  '
